# Rename "Uncut_Sheet" -> "Uncut_Sheet_1", make it the active/selected tab
# (it previously was not active; "Pipette" was), and change its selection
# from D17:E17 to B17:C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Uncut_Sheet")

# Rename the sheet. (Single-quoted literal below avoids PowerShell treating
# "$A$1" etc. as variable references.)
$ws.Name = "Uncut_Sheet_1"

# Re-assert the sheet's print area so the workbook-level
# _xlnm.Print_Area defined name is rewritten to point at the new sheet
# name instead of the stale "Uncut_Sheet" reference.
$ws.PageSetup.PrintArea = '$A$1:$G$42'

# Switch the active tab to this sheet and move its selection to B17:C17.
$ws.Activate() | Out-Null
$ws.Range("B17:C17").Select() | Out-Null
